$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushing existing rows 3-11 down to 4-12),
# carrying formatting (fill colors etc.) down with the shift - this
# matches how Excel's native "Insert Row" behaves.
$ws.Rows.Item(3).Insert()

# Fill the newly inserted row 3 with the "9:00 - 10:00" interval label
# and the same green fill used by the other interval rows (copied from
# row 2, which already has the correct style for columns B:H).
$ws.Range("A3").Value = "9:00 - 10:00"
$ws.Range("B2:H2").Copy() | Out-Null
$ws.Range("B3:H3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the active selection to match the post-edit state.
$ws.Range("A3").Select() | Out-Null
